$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.834.09"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "1.876.03"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "'0.7208"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.33%  "
$ws.Range("D6").Value = "'242.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").Value = "'0.3150"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("D9").Value = "'0.07369"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.26%  "
$ws.Range("D10").Value = "'24.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.91%  "
$ws.Range("D11").Value = "'0.08202"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.38%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.900.77"
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.7452"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").Value = "'5.330"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").Value = "'92.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").Value = "29.789.57"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").Value = "'6.014"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("D18").Value = "'246.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").Value = "'13.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("D20").Value = "'0.000007897"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.97%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").Value = "2.116.74"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "'7.727"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.14%  "
$ws.Range("D25").Value = "'9.243"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("D26").Value = "'0.1505"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.41%  "
$ws.Range("D27").Value = "'164.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").Value = "'18.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("D29").Value = "'2.009"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("D30").Value = "'1.423"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.74%  "
$ws.Range("E31").Value = "  -1.38%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "'4.170"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.37%  "
$ws.Range("D34").Value = "'0.05482"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.07%  "
$ws.Range("D35").Value = "'1.230"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("D36").Value = "'0.7340"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.71%  "
$ws.Range("D37").Value = "'0.9995"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").Value = "'2.703"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").Value = "'0.01913"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").Value = "'2.741"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.49%  "
$ws.Range("D41").Value = "'0.4450"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.98%  "
$ws.Range("D42").Value = "'0.8937"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.41%  "
$ws.Range("D43").Value = "'5.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").Value = "'71.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").Value = "1.039.02"
$ws.Range("E46").Value = "  -6.61%  "
$ws.Range("D47").Value = "'103.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("D48").Value = "'7.469"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.37%  "
$ws.Range("D49").Value = "'1.810"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.74%  "
$ws.Range("D50").Value = "'9.607"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("D51").Value = "2.016.67"
$ws.Range("E51").Value = "  -0.17%  "
